$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix the mis-keyed date on row 75 (was 41382, should be 45034)
$ws.Range("A75").Value = 45034

# Copy row 75's formatting (date / time number formats) down to the two
# new rows before filling in their values, so the new cells pick up the
# same cell styles (s="1" on col A, s="2" on col B) instead of creating
# brand-new style entries.
$ws.Range("A75:D75").Copy() | Out-Null
$ws.Range("A76:D76").PasteSpecial(-4122) | Out-Null
$ws.Range("A77:D77").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# New data rows scraped from supremacy1914 / disk_savvy logs
$ws.Range("A76").Value = 45088
$ws.Range("B76").Value = 0.44122685185185184
$ws.Range("C76").Value = 77494
$ws.Range("D76").Value = 1480

$ws.Range("A77").Value = 45091
$ws.Range("B77").Value = 0.44196759259259261
$ws.Range("C77").Value = 81206
$ws.Range("D77").Value = 1620

# Move the selection to reflect where the user ended up after typing the
# new rows (Excel leaves the cursor one row below the last entry).
$ws.Range("A78").Select() | Out-Null
